# Fixed Feign Relations Finder and added 7 projects.
#
# This script re-applies the "classFields" sheet field-name/field-type
# pairings so they line up correctly (the relations finder previously
# associated some field names with the wrong declared type). Only the
# B (Field Name) and D (Field Type) columns move; A (Class Name) and
# C (Field Modifier) stay put for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# Row 2-4 (com.macro.mall.config.SwaggerResourceConfig fields) cycle:
#   log/Logger -> routeLocator/RouteLocator -> gatewayProperties/GatewayProperties -> log/Logger
$ws.Range("B2").Value = "gatewayProperties"
$ws.Range("D2").Value = "org.springframework.cloud.gateway.config.GatewayProperties"

$ws.Range("B3").Value = "log"
$ws.Range("D3").Value = "org.slf4j.Logger"

$ws.Range("B4").Value = "routeLocator"
$ws.Range("D4").Value = "org.springframework.cloud.gateway.route.RouteLocator"

# Row 6-8 (com.macro.mall.handler.SwaggerHandler fields) cycle:
#   uiConfiguration -> securityConfiguration -> swaggerResources -> uiConfiguration
$ws.Range("B6").Value = "securityConfiguration"
$ws.Range("D6").Value = "springfox.documentation.swagger.web.SecurityConfiguration"

$ws.Range("B7").Value = "swaggerResources"
$ws.Range("D7").Value = "springfox.documentation.swagger.web.SwaggerResourcesProvider"

$ws.Range("B8").Value = "uiConfiguration"
$ws.Range("D8").Value = "springfox.documentation.swagger.web.UiConfiguration"

# Rows 10,12,13,14 (com.macro.mall.config.ResourceServerConfig fields) cycle:
#   ignoreUrlsConfig(10) -> ignoreUrlsRemoveJwtFilter
#   ignoreUrlsRemoveJwtFilter(12) -> restAuthenticationEntryPoint
#   restfulAccessDeniedHandler(13) -> ignoreUrlsConfig
#   restAuthenticationEntryPoint(14) -> restfulAccessDeniedHandler
$ws.Range("B10").Value = "ignoreUrlsRemoveJwtFilter"
$ws.Range("D10").Value = "com.macro.mall.filter.IgnoreUrlsRemoveJwtFilter"

$ws.Range("B12").Value = "restAuthenticationEntryPoint"
$ws.Range("D12").Value = "com.macro.mall.component.RestAuthenticationEntryPoint"

$ws.Range("B13").Value = "ignoreUrlsConfig"
$ws.Range("D13").Value = "com.macro.mall.config.IgnoreUrlsConfig"

$ws.Range("B14").Value = "restfulAccessDeniedHandler"
$ws.Range("D14").Value = "com.macro.mall.component.RestfulAccessDeniedHandler"
